# Auto-generated edit script applying the Carbuncle_Profits (workbook) diff.
# Updates currentAveragePrice / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns (H,I,J,K,L,M,N) for the affected leve rows
# on each class sheet, matching a scheduled market-price refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1319.3572
$ws.Range("J32").Value = 1422.5834
$ws.Range("L32").Value = 1422.5834
$ws.Range("N32").Value = -2074.5834

# Row 88
$ws.Range("H88").Value = 2153.8
$ws.Range("J88").Value = 2193.3845
$ws.Range("L88").Value = 2193.3845
$ws.Range("N88").Value = -3005.3845

# Row 91
$ws.Range("H91").Value = 2153.8
$ws.Range("J91").Value = 2193.3845
$ws.Range("L91").Value = 2193.3845
$ws.Range("N91").Value = -5001.3845

# Row 111
$ws.Range("H111").Value = 100002340
$ws.Range("I111").Value = 1929.5
$ws.Range("K111").Value = 5788.5
$ws.Range("M111").Value = -2721.5

# Row 113
$ws.Range("H113").Value = 8044.3477
$ws.Range("I113").Value = 2309.2307
$ws.Range("K113").Value = 2309.2307
$ws.Range("M113").Value = 944.7692999999999

# Row 116
$ws.Range("H116").Value = 3128.4285
$ws.Range("I116").Value = 3299.75
$ws.Range("J116").Value = 2900
$ws.Range("K116").Value = 3299.75
$ws.Range("L116").Value = 2900
$ws.Range("M116").Value = 142.25
$ws.Range("N116").Value = -9784

# Row 132
$ws.Range("H132").Value = 4955.9565
$ws.Range("I132").Value = 5623.125
$ws.Range("J132").Value = 3431
$ws.Range("K132").Value = 16869.375
$ws.Range("L132").Value = 10293
$ws.Range("M132").Value = -14339.375
$ws.Range("N132").Value = -15353

# Row 137
$ws.Range("H137").Value = 1666.3226
$ws.Range("I137").Value = 1324.6666
$ws.Range("J137").Value = 2139.3845
$ws.Range("K137").Value = 3973.9998
$ws.Range("L137").Value = 6418.1535
$ws.Range("M137").Value = -1423.9998
$ws.Range("N137").Value = -11518.1535

# Row 141
$ws.Range("H141").Value = 3515
$ws.Range("I141").Value = 3515
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10545
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5365
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 3549.4211
$ws.Range("I63").Value = 3643.4707
$ws.Range("K63").Value = 3643.4707
$ws.Range("M63").Value = -2957.4707

# Row 66
$ws.Range("H66").Value = 3549.4211
$ws.Range("I66").Value = 3643.4707
$ws.Range("K66").Value = 18217.3535
$ws.Range("M66").Value = -14785.3535

# Row 109
$ws.Range("H109").Value = 45377
$ws.Range("J109").Value = 45377
$ws.Range("L109").Value = 45377
$ws.Range("N109").Value = -48151

# Row 110
$ws.Range("H110").Value = 22215.75
$ws.Range("I110").Value = 30857.824
$ws.Range("J110").Value = 1227.8572
$ws.Range("K110").Value = 30857.824
$ws.Range("L110").Value = 1227.8572
$ws.Range("M110").Value = -28812.824
$ws.Range("N110").Value = -5317.8572

# Row 115
$ws.Range("H115").Value = 40184
$ws.Range("J115").Value = 40184
$ws.Range("L115").Value = 40184
$ws.Range("N115").Value = -43318

# Row 132
$ws.Range("H132").Value = 2379.7896
$ws.Range("I132").Value = 1800.2858
$ws.Range("J132").Value = 4002.4
$ws.Range("K132").Value = 5400.857400000001
$ws.Range("L132").Value = 12007.2
$ws.Range("M132").Value = -2870.857400000001
$ws.Range("N132").Value = -17067.2

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2161.95
$ws.Range("I134").Value = 2248.4285
$ws.Range("J134").Value = 2115.3845
$ws.Range("K134").Value = 6745.2855
$ws.Range("L134").Value = 6346.1535
$ws.Range("M134").Value = -4210.2855
$ws.Range("N134").Value = -11416.1535

# Row 140
$ws.Range("H140").Value = 42492
$ws.Range("J140").Value = 42492
$ws.Range("L140").Value = 42492
$ws.Range("N140").Value = -52852

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 11832
$ws.Range("I99").Value = 2083.6667
$ws.Range("K99").Value = 2083.6667
$ws.Range("M99").Value = -585.6667000000002

# Row 126
$ws.Range("H126").Value = 11832
$ws.Range("I126").Value = 2083.6667
$ws.Range("K126").Value = 6251.000100000001
$ws.Range("M126").Value = -3781.000100000001

# Row 132
$ws.Range("H132").Value = 2992.8696
$ws.Range("I132").Value = 2021.6
$ws.Range("J132").Value = 4814
$ws.Range("K132").Value = 6064.799999999999
$ws.Range("L132").Value = 14442
$ws.Range("M132").Value = -3534.799999999999
$ws.Range("N132").Value = -19502

# Row 135
$ws.Range("H135").Value = 43000
$ws.Range("J135").Value = 43000
$ws.Range("L135").Value = 43000
$ws.Range("N135").Value = -53140

# Row 140
$ws.Range("H140").Value = 70804
$ws.Range("J140").Value = 70804
$ws.Range("L140").Value = 70804
$ws.Range("N140").Value = -81164

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1377.4535
$ws.Range("I131").Value = 614.2308
$ws.Range("J131").Value = 1513.3699
$ws.Range("K131").Value = 1842.6924
$ws.Range("L131").Value = 4540.1097
$ws.Range("M131").Value = 3197.3076
$ws.Range("N131").Value = -14620.1097

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 19800
$ws.Range("J93").Value = 19800
$ws.Range("L93").Value = 19800
$ws.Range("N93").Value = -23544

# Row 126
$ws.Range("H126").Value = 2205.3914
$ws.Range("I126").Value = 1889
$ws.Range("J126").Value = 2550.5454
$ws.Range("K126").Value = 5667
$ws.Range("L126").Value = 7651.6362
$ws.Range("M126").Value = -3197
$ws.Range("N126").Value = -12591.6362

# Row 138
$ws.Range("H138").Value = 45722.6
$ws.Range("J138").Value = 45722.6
$ws.Range("L138").Value = 45722.6
$ws.Range("N138").Value = -56002.6

# Row 140
$ws.Range("H140").Value = 39044464
$ws.Range("J140").Value = 39044464
$ws.Range("L140").Value = 39044464
$ws.Range("N140").Value = -39054824

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 306.8421
$ws.Range("I55").Value = 228.33333
$ws.Range("J55").Value = 441.42856
$ws.Range("K55").Value = 228.33333
$ws.Range("L55").Value = 441.42856
$ws.Range("M55").Value = -55.33332999999999
$ws.Range("N55").Value = -787.4285600000001

# Row 122
$ws.Range("H122").Value = 2294
$ws.Range("I122").Value = 2140
$ws.Range("J122").Value = 2448
$ws.Range("K122").Value = 6420
$ws.Range("L122").Value = 7344
$ws.Range("M122").Value = -3970
$ws.Range("N122").Value = -12244

# Row 132
$ws.Range("H132").Value = 10422.206
$ws.Range("I132").Value = 11682.913
$ws.Range("J132").Value = 7786.1816
$ws.Range("K132").Value = 35048.739
$ws.Range("L132").Value = 23358.5448
$ws.Range("M132").Value = -32518.739
$ws.Range("N132").Value = -28418.5448

$ws = $wb.Worksheets.Item("WVR")
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

# Row 137
$ws.Range("H137").Value = 67270.25
$ws.Range("J137").Value = 67270.25
$ws.Range("L137").Value = 67270.25
$ws.Range("N137").Value = -77470.25
